$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string value into a cell while forcing Excel to keep it
# as text (avoids "39.460.99", "0.999", "16.01", etc. being auto-converted
# into numbers / losing trailing zeros). ClearFormats() afterwards removes
# the temporary text NumberFormat again so the cells style is left exactly
# as it was originally (no explicit style / "General" format).
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "39.460.99"
Set-TextValue $ws.Range("E2") "  +1.76%  "

Set-TextValue $ws.Range("D3") "2.158.81"
Set-TextValue $ws.Range("E3") "  +2.73%  "

Set-TextValue $ws.Range("E4") "  +0.06%  "

Set-TextValue $ws.Range("D5") "227.90"
Set-TextValue $ws.Range("E5") "  -0.52%  "

Set-TextValue $ws.Range("E6") "  +0.85%  "

Set-TextValue $ws.Range("D7") "64.10"
Set-TextValue $ws.Range("E7") "  +4.01%  "

Set-TextValue $ws.Range("E8") "  +0.05%  "

Set-TextValue $ws.Range("E9") "  +2.57%  "

Set-TextValue $ws.Range("D10") "0.0856"
Set-TextValue $ws.Range("E10") "  +1.22%  "

Set-TextValue $ws.Range("E11") "  +0.23%  "

Set-TextValue $ws.Range("D12") "16.01"
Set-TextValue $ws.Range("E12") "  +3.18%  "

Set-TextValue $ws.Range("D13") "2.479.31"
Set-TextValue $ws.Range("E13") "  +2.77%  "

Set-TextValue $ws.Range("E14") "  +0.56%  "

Set-TextValue $ws.Range("E15") "  +0.40%  "

Set-TextValue $ws.Range("E16") "  +0.71%  "

Set-TextValue $ws.Range("D17") "2.149.39"
Set-TextValue $ws.Range("E17") "  +2.17%  "

Set-TextValue $ws.Range("D18") "39.413.23"
Set-TextValue $ws.Range("E18") "  +1.55%  "

Set-TextValue $ws.Range("D19") "71.84"
Set-TextValue $ws.Range("E19") "  -0.24%  "

Set-TextValue $ws.Range("E20") "  +0.42%  "

Set-TextValue $ws.Range("D21") "0.0₃0850"
Set-TextValue $ws.Range("E21") "  +1.19%  "

Set-TextValue $ws.Range("D22") "231.16"
Set-TextValue $ws.Range("E22") "  +1.56%  "

Set-TextValue $ws.Range("E23") "  +0.01%  "

Set-TextValue $ws.Range("E24") "  +0.47%  "

Set-TextValue $ws.Range("D25") "2.31"
Set-TextValue $ws.Range("E25") "  -3.01%  "

Set-TextValue $ws.Range("D26") "172.24"
Set-TextValue $ws.Range("E26") "  +0.40%  "

Set-TextValue $ws.Range("E27") "  -0.44%  "

Set-TextValue $ws.Range("E28") "  +1.09%  "

Set-TextValue $ws.Range("E29") "  +2.68%  "

Set-TextValue $ws.Range("E30") "  -0.45%  "

Set-TextValue $ws.Range("D31") "2.66"
Set-TextValue $ws.Range("E31") "  +7.70%  "

Set-TextValue $ws.Range("E32") "  +0.54%  "

Set-TextValue $ws.Range("D33") "4.61"
Set-TextValue $ws.Range("E33") "  +1.74%  "

Set-TextValue $ws.Range("E34") "  +9.58%  "

Set-TextValue $ws.Range("D35") "4.74"
Set-TextValue $ws.Range("E35") "  -0.67%  "

Set-TextValue $ws.Range("E36") "  -0.77%  "

Set-TextValue $ws.Range("E37") "  +0.32%  "

Set-TextValue $ws.Range("D38") "3.58"
Set-TextValue $ws.Range("E38") "  +0.14%  "

Set-TextValue $ws.Range("D39") "0.999"
Set-TextValue $ws.Range("E39") "  -0.12%  "

Set-TextValue $ws.Range("D40") "103.68"
Set-TextValue $ws.Range("E40") "  +2.12%  "

Set-TextValue $ws.Range("E41") "  +0.83%  "

Set-TextValue $ws.Range("D42") "17.69"
Set-TextValue $ws.Range("E42") "  -2.61%  "

Set-TextValue $ws.Range("D43") "1.540.58"
Set-TextValue $ws.Range("E43") "  +0.44%  "

Set-TextValue $ws.Range("E44") "  +4.33%  "

Set-TextValue $ws.Range("E45") "  +5.96%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D46") "0.0931"
Set-TextValue $ws.Range("E46") "  +2.26%  "

$ws.Range("B47").Value = "HuobiToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D47") "2.83"
Set-TextValue $ws.Range("E47") "  +0.82%  "

Set-TextValue $ws.Range("E48") "  +4.98%  "

Set-TextValue $ws.Range("D49") "7.70"
Set-TextValue $ws.Range("E49") "  -0.49%  "

Set-TextValue $ws.Range("D50") "2.363.15"
Set-TextValue $ws.Range("E50") "  +2.96%  "

Set-TextValue $ws.Range("E51") "  -0.12%  "
